$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1 with same style as the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Fill H2:H18 with 0 for each data row
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
